$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-01 Thursday" "2026-01-02 Friday"

Replace-Text "273÷7=" "285÷7="
Replace-Text "144÷2=" "373÷6="
Replace-Text "653÷5=" "265÷9="
Replace-Text "554÷4=" "345÷7="
Replace-Text "530÷7=" "562÷8="
Replace-Text "620÷7=" "908÷9="
Replace-Text "790÷4=" "678÷5="
Replace-Text "100÷2=" "751÷7="
Replace-Text "414÷5=" "540÷4="
Replace-Text "639÷6=" "307÷6="
Replace-Text "374÷2=" "989÷5="
Replace-Text "232÷5=" "949÷4="
Replace-Text "123÷2=" "617÷7="
Replace-Text "270÷4=" "863÷6="
Replace-Text "115÷3=" "619÷3="
Replace-Text "749÷9=" "890÷5="
Replace-Text "847÷3=" "669÷7="
Replace-Text "791÷3=" "300÷5="
Replace-Text "399÷8=" "476÷6="
Replace-Text "937÷7=" "360÷7="
Replace-Text "140÷2=" "587÷9="
Replace-Text "462÷6=" "677÷3="
Replace-Text "971÷6=" "181÷5="
Replace-Text "776÷9=" "807÷9="
Replace-Text "735÷3=" "940÷9="
